$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Prn"
$ws.Cells.Item(2, 3).Value = "Rpsa"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.7438956666666666
$ws.Cells.Item(2, 8).Value = 2.231687
$ws.Cells.Item(2, 9).Value = 0.2947327619223709
$ws.Cells.Item(2, 10).Value = 0.2947327619223709
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 112.708133
$ws.Cells.Item(2, 14).Value = 338.124399
$ws.Cells.Item(2, 15).Value = 0.2121524692929861
$ws.Cells.Item(2, 16).Value = 0.2121524692929861
$ws.Cells.Item(2, 17).Value = 83.84309173679034
$ws.Cells.Item(2, 18).Value = 754.5878256311131
$ws.Cells.Item(2, 19).Value = 0.06252828322337277
$ws.Cells.Item(2, 20).Value = 0.06252828322337277

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Prn"
$ws.Cells.Item(3, 3).Value = "Rpsa"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.7438956666666666
$ws.Cells.Item(3, 8).Value = 2.231687
$ws.Cells.Item(3, 9).Value = 0.2947327619223709
$ws.Cells.Item(3, 10).Value = 0.2947327619223709
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 189.57842
$ws.Cells.Item(3, 14).Value = 568.7352599999999
$ws.Cells.Item(3, 15).Value = 0.3568467408440064
$ws.Cells.Item(3, 16).Value = 0.3568467408440064
$ws.Cells.Item(3, 17).Value = 141.0265651315133
$ws.Cells.Item(3, 18).Value = 1269.23908618362
$ws.Cells.Item(3, 19).Value = 0.1051744255119505
$ws.Cells.Item(3, 20).Value = 0.1051744255119505

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Prn"
$ws.Cells.Item(4, 3).Value = "Rpsa"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.7438956666666666
$ws.Cells.Item(4, 8).Value = 2.231687
$ws.Cells.Item(4, 9).Value = 0.2947327619223709
$ws.Cells.Item(4, 10).Value = 0.2947327619223709
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 110.6512236666667
$ws.Cells.Item(4, 14).Value = 331.953671
$ws.Cells.Item(4, 15).Value = 0.2082807132576123
$ws.Cells.Item(4, 16).Value = 0.2082807132576123
$ws.Cells.Item(4, 17).Value = 82.31296579699745
$ws.Cells.Item(4, 18).Value = 740.816692172977
$ws.Cells.Item(4, 19).Value = 0.06138714987357743
$ws.Cells.Item(4, 20).Value = 0.06138714987357744

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Prn"
$ws.Cells.Item(5, 3).Value = "Rpsa"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.7438956666666666
$ws.Cells.Item(5, 8).Value = 2.231687
$ws.Cells.Item(5, 9).Value = 0.2947327619223709
$ws.Cells.Item(5, 10).Value = 0.2947327619223709
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 118.3222806666667
$ws.Cells.Item(5, 14).Value = 354.966842
$ws.Cells.Item(5, 15).Value = 0.2227200766053953
$ws.Cells.Item(5, 16).Value = 0.2227200766053952
$ws.Cells.Item(5, 17).Value = 88.01943185805044
$ws.Cells.Item(5, 18).Value = 792.174886722454
$ws.Cells.Item(5, 19).Value = 0.06564290331347016
$ws.Cells.Item(5, 20).Value = 0.06564290331347016

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Prn"
$ws.Cells.Item(6, 3).Value = "Rpsa"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.780071
$ws.Cells.Item(6, 8).Value = 5.340212999999999
$ws.Cells.Item(6, 9).Value = 0.705267238077629
$ws.Cells.Item(6, 10).Value = 0.7052672380776291
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 112.708133
$ws.Cells.Item(6, 14).Value = 338.124399
$ws.Cells.Item(6, 15).Value = 0.2121524692929861
$ws.Cells.Item(6, 16).Value = 0.2121524692929861
$ws.Cells.Item(6, 17).Value = 200.628479017443
$ws.Cells.Item(6, 18).Value = 1805.656311156987
$ws.Cells.Item(6, 19).Value = 0.1496241860696133
$ws.Cells.Item(6, 20).Value = 0.1496241860696133

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Prn"
$ws.Cells.Item(7, 3).Value = "Rpsa"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.780071
$ws.Cells.Item(7, 8).Value = 5.340212999999999
$ws.Cells.Item(7, 9).Value = 0.705267238077629
$ws.Cells.Item(7, 10).Value = 0.7052672380776291
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 189.57842
$ws.Cells.Item(7, 14).Value = 568.7352599999999
$ws.Cells.Item(7, 15).Value = 0.3568467408440064
$ws.Cells.Item(7, 16).Value = 0.3568467408440064
$ws.Cells.Item(7, 17).Value = 337.4630476678199
$ws.Cells.Item(7, 18).Value = 3037.167429010379
$ws.Cells.Item(7, 19).Value = 0.2516723153320559
$ws.Cells.Item(7, 20).Value = 0.2516723153320559

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Prn"
$ws.Cells.Item(8, 3).Value = "Rpsa"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.780071
$ws.Cells.Item(8, 8).Value = 5.340212999999999
$ws.Cells.Item(8, 9).Value = 0.705267238077629
$ws.Cells.Item(8, 10).Value = 0.7052672380776291
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 110.6512236666667
$ws.Cells.Item(8, 14).Value = 331.953671
$ws.Cells.Item(8, 15).Value = 0.2082807132576123
$ws.Cells.Item(8, 16).Value = 0.2082807132576123
$ws.Cells.Item(8, 17).Value = 196.967034363547
$ws.Cells.Item(8, 18).Value = 1772.703309271923
$ws.Cells.Item(8, 19).Value = 0.1468935633840348
$ws.Cells.Item(8, 20).Value = 0.1468935633840348

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Prn"
$ws.Cells.Item(9, 3).Value = "Rpsa"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.780071
$ws.Cells.Item(9, 8).Value = 5.340212999999999
$ws.Cells.Item(9, 9).Value = 0.705267238077629
$ws.Cells.Item(9, 10).Value = 0.7052672380776291
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 118.3222806666667
$ws.Cells.Item(9, 14).Value = 354.966842
$ws.Cells.Item(9, 15).Value = 0.2227200766053953
$ws.Cells.Item(9, 16).Value = 0.2227200766053952
$ws.Cells.Item(9, 17).Value = 210.622060468594
$ws.Cells.Item(9, 18).Value = 1895.598544217346
$ws.Cells.Item(9, 19).Value = 0.1570771732919251
$ws.Cells.Item(9, 20).Value = 0.1570771732919251
